$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''250.32'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = '''22.71'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = '''5.420'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = '''0.05730'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = '''3.413'
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = '''6.343'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = '''0.8137'
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = '''0.9351'
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = '''0.1424'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = '''0.07554'
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = '''0.03132'
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = '''0.03085'
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = '''0.09365'
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = '''3.735'
$ws.Range("D15").Style = "Normal"

# Row 17
$ws.Range("D17").Value = '''0.04771'
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = '''UpBots'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '''https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '''0.01829'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''17UpBotsUBXTBestin24h'
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = '''One'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '''https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '''0.0005795'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''18OneONE'
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = '''TigerCash'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '''0.006471'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''19TigerCashTCH'
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("B21").Value = '''HotbitToken'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '''https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '''0.005010'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''20HotbitTokenHTB'
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("B22").Value = '''BitKan'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '''https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '''0.001026'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''21BitKanKAN'
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("B23").Value = '''NitroEx'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '''https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = '''0.0001500'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''22NitroExNTX'
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("B24").Value = '''LEO'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = '''3.702'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''23LEOLEO'
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("B25").Value = '''BTSEToken'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '''https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = '''2.163'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''24BTSETokenBTSE'
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = '''BitpandaEcosystemToken'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '''https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = '''0.3303'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''25BitpandaEcosystemTokenBEST'
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = '''ProBitToken'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = '''https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = '''0.1309'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''26ProBitTokenPROB'
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = '''AAXToken'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = '''https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = '''0.1950'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''27AAXTokenAAB'
$ws.Range("E28").Style = "Normal"

# Row 40
$ws.Range("D40").Value = '''0.04011'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = '''BKEXToken'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.1069'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''40BKEXTokenBKK'
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = '''CEJI'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''0.002710'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''41CEJICEJI'
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = '''KickToken'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''0.002945'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''42KickTokenKICKWorstin24h'
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = '''0.007956'
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = '''0.00005895'
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = '''0.00000000751'
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = '''0.5005'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''46CoinbaseStockTokenCOIN'
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = '''47BOLOBOLO'
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = '''0.00002102'
$ws.Range("D49").Style = "Normal"
